$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as Text so numeric-looking values
# (e.g. "1.002") are preserved as strings rather than being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.494.44"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "1.637.19"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "306.88"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.3762"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "52.60"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "0.3643"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "1.269"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "0.08174"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "22.97"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "6.633"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "7.382"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "1.637.53"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "94.74"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "0.06956"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "18.23"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "6.553"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "23.518.61"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "12.81"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "3.084"
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("D26").Value = "2.420"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "21.28"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "151.34"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "5.356"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("D30").Value = "135.49"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").Value = "2.366"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "1.819.59"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "6.790"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").Value = "0.9648"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "0.02819"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "0.07350"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").Value = "0.2538"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "6.179"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "0.08861"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").Value = "0.7110"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "16.25"
$ws.Range("E43").Value = "  +5.69%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "12.49"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "0.6548"
$ws.Range("D46").Value = "2.343"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "4.031"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "0.07969"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "129.43"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").Value = "1.207"
$ws.Range("E51").Value = "  +0.26%  "
